$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they stay text (matching original inlineStr type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range('D2').Value = '61.032.89'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '3.391.63'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '572.62'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '142.10'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('D7').Value = '3.392.13'
$ws.Range('E7').Value = '  -1.53%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').Value = '7.64'
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('D13').Value = '3.971.61'
$ws.Range('E13').Value = '  -1.56%  '
$ws.Range('D14').Value = '0.125'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '27.89'
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('D17').Value = '3.377.43'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '61.081.43'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('D20').Value = '13.77'
$ws.Range('E20').Value = '  -4.73%  '
$ws.Range('E21').Value = '  -5.10%  '
$ws.Range('D22').Value = '382.08'
$ws.Range('E22').Value = '  -4.18%  '
$ws.Range('D23').Value = '0.555'
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').Value = '74.55'
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('E26').Value = '  -4.78%  '
$ws.Range('D27').Value = '3.530.80'
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = '7.33'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.15'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '7.97'
$ws.Range('E32').Value = '  -3.24%  '
$ws.Range('E33').Value = '  -6.03%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '23.44'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('D37').Value = '166.43'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('D38').Value = '3.424.49'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('E41').Value = '  -2.55%  '
$ws.Range('D42').Value = '27.24'
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '42.12'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '4.41'
$ws.Range('E46').Value = '  -2.62%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '1.67'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '1.14'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.475.11'
$ws.Range('E49').Value = '  -5.52%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '6.79'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '22.86'
$ws.Range('E51').Value = '  -1.10%  '
